$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Actual Risk ($)" column entirely (it duplicated "Risk Amount ($)")
$ws.Range("N1").EntireColumn.Delete()

# Re-apply a flat currency format (no bold/fill) to the Entry Price and PnL ($) columns,
# which previously used the separate "Currency" accounting style.
$ws.Range("E1:E37").Style = "Normal"
$ws.Range("E2:E37").NumberFormat = "\$#,##0.00"
$ws.Range("K2:K37").NumberFormat = "\$#,##0.00"

# Move the active selection
$ws.Range("C42").Select()
